$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new sheet "Tasks 02-04 to 02-11" between "TaskList" and
#    "Tasks 01-28 to 02-04". Worksheets.Add() inserts right before the
#    currently active sheet, which is "Tasks 01-28 to 02-04".
# ------------------------------------------------------------------
$wb.Worksheets.Item("Tasks 01-28 to 02-04").Activate()
$new = $wb.Worksheets.Add()
$new.Name = "Tasks 02-04 to 02-11"

# NOTE: worksheet variables are resolved positionally, so any reference
# captured before the Add() call above now points at the *new* sheet.
# Re-fetch the original sheet by name now that it has moved to index 3.
$old = $wb.Worksheets.Item("Tasks 01-28 to 02-04")

# Cells already carrying the exact styles we need to reuse, taken from
# the original "Tasks 01-28 to 02-04" sheet (copied by reference so we
# never create duplicate font/fill entries in styles.xml).
$style25 = $old.Range("A1")
$style14 = $old.Range("A3")
$style26 = $old.Range("M9")
$style27 = $old.Range("F8")
$style28 = $old.Range("F2")

# ---- Row values -----------------------------------------------------
# Row 1
$new.Range("A1").Value = "Task"
$style25.Copy() | Out-Null
$new.Range("A1").PasteSpecial(-4122) | Out-Null
$new.Range("B1").Value = "Time Estimated to Complete"
$style25.Copy() | Out-Null
$new.Range("B1").PasteSpecial(-4122) | Out-Null
$new.Range("C1").Value = "Time Spent"
$style25.Copy() | Out-Null
$new.Range("C1").PasteSpecial(-4122) | Out-Null
$new.Range("D1").Value = "Over/Under"
$style25.Copy() | Out-Null
$new.Range("D1").PasteSpecial(-4122) | Out-Null
$new.Range("E1").Value = "Assigned To"
$style25.Copy() | Out-Null
$new.Range("E1").PasteSpecial(-4122) | Out-Null
$new.Range("F1").Value = "Status"
$style25.Copy() | Out-Null
$new.Range("F1").PasteSpecial(-4122) | Out-Null
$new.Range("G1").Value = "Notes"
$style25.Copy() | Out-Null
$new.Range("G1").PasteSpecial(-4122) | Out-Null

# Row 2
$new.Range("A2").Value = "Add a RigidBody to the base icon"
$style25.Copy() | Out-Null
$new.Range("A2").PasteSpecial(-4122) | Out-Null
$new.Range("B2").Value = 0.5
$style25.Copy() | Out-Null
$new.Range("B2").PasteSpecial(-4122) | Out-Null
$style25.Copy() | Out-Null
$new.Range("C2").PasteSpecial(-4122) | Out-Null
$style25.Copy() | Out-Null
$new.Range("D2").PasteSpecial(-4122) | Out-Null
$new.Range("E2").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E2").PasteSpecial(-4122) | Out-Null
$new.Range("F2").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F2").PasteSpecial(-4122) | Out-Null
$style25.Copy() | Out-Null
$new.Range("G2").PasteSpecial(-4122) | Out-Null

# Row 3
$new.Range("A3").Value = "Create MagicIcon"
$style14.Copy() | Out-Null
$new.Range("A3").PasteSpecial(-4122) | Out-Null
$new.Range("B3").Value = 2
$style14.Copy() | Out-Null
$new.Range("B3").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C3").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D3").PasteSpecial(-4122) | Out-Null
$new.Range("E3").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E3").PasteSpecial(-4122) | Out-Null
$new.Range("F3").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F3").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("G3").PasteSpecial(-4122) | Out-Null

# Row 4
$new.Range("A4").Value = "Create CoinIcon"
$style14.Copy() | Out-Null
$new.Range("A4").PasteSpecial(-4122) | Out-Null
$new.Range("B4").Value = 2
$style14.Copy() | Out-Null
$new.Range("B4").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C4").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D4").PasteSpecial(-4122) | Out-Null
$new.Range("E4").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E4").PasteSpecial(-4122) | Out-Null
$new.Range("F4").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F4").PasteSpecial(-4122) | Out-Null
$new.Range("G4").Value = "Adds a coin when colliding with coinbag"
$style14.Copy() | Out-Null
$new.Range("G4").PasteSpecial(-4122) | Out-Null

# Row 5
$new.Range("A5").Value = "Create RangedIcon"
$style14.Copy() | Out-Null
$new.Range("A5").PasteSpecial(-4122) | Out-Null
$new.Range("B5").Value = 2
$style14.Copy() | Out-Null
$new.Range("B5").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C5").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D5").PasteSpecial(-4122) | Out-Null
$new.Range("E5").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E5").PasteSpecial(-4122) | Out-Null
$new.Range("F5").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F5").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("G5").PasteSpecial(-4122) | Out-Null

# Row 6
$new.Range("A6").Value = "Have IconSpawner Spawn in new icons when old ones are used"
$style14.Copy() | Out-Null
$new.Range("A6").PasteSpecial(-4122) | Out-Null
$new.Range("B6").Value = 2.5
$style14.Copy() | Out-Null
$new.Range("B6").PasteSpecial(-4122) | Out-Null
$new.Range("E6").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E6").PasteSpecial(-4122) | Out-Null
$new.Range("F6").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F6").PasteSpecial(-4122) | Out-Null

# Row 7
$new.Range("A7").Value = "Create Base Enemy Class"
$style14.Copy() | Out-Null
$new.Range("A7").PasteSpecial(-4122) | Out-Null
$new.Range("B7").Value = 3
$style14.Copy() | Out-Null
$new.Range("B7").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C7").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D7").PasteSpecial(-4122) | Out-Null
$new.Range("E7").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E7").PasteSpecial(-4122) | Out-Null
$new.Range("F7").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F7").PasteSpecial(-4122) | Out-Null
$new.Range("G7").Value = "Base Enemy class will have basic enemy behaviours"
$style14.Copy() | Out-Null
$new.Range("G7").PasteSpecial(-4122) | Out-Null

# Row 8
$new.Range("A8").Value = "Create EnemyPrefab"
$style14.Copy() | Out-Null
$new.Range("A8").PasteSpecial(-4122) | Out-Null
$new.Range("B8").Value = 0.5
$style14.Copy() | Out-Null
$new.Range("B8").PasteSpecial(-4122) | Out-Null
$new.Range("E8").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E8").PasteSpecial(-4122) | Out-Null
$new.Range("F8").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F8").PasteSpecial(-4122) | Out-Null
$new.Range("H8").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("H8").PasteSpecial(-4122) | Out-Null

# Row 9
$new.Range("A9").Value = "Create EnemySpawner"
$style14.Copy() | Out-Null
$new.Range("A9").PasteSpecial(-4122) | Out-Null
$new.Range("B9").Value = 1
$style14.Copy() | Out-Null
$new.Range("B9").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C9").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D9").PasteSpecial(-4122) | Out-Null
$new.Range("E9").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E9").PasteSpecial(-4122) | Out-Null
$new.Range("F9").Value = "TODO"
$style26.Copy() | Out-Null
$new.Range("F9").PasteSpecial(-4122) | Out-Null
$new.Range("G9").Value = "Will Spawn in Enemy Prefab"
$style14.Copy() | Out-Null
$new.Range("G9").PasteSpecial(-4122) | Out-Null
$new.Range("H9").Value = "In Progress"
$style27.Copy() | Out-Null
$new.Range("H9").PasteSpecial(-4122) | Out-Null

# Row 10
$new.Range("A10").Value = "Create ComparedItem Game Object"
$style14.Copy() | Out-Null
$new.Range("A10").PasteSpecial(-4122) | Out-Null
$new.Range("B10").Value = 2
$style14.Copy() | Out-Null
$new.Range("B10").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C10").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D10").PasteSpecial(-4122) | Out-Null
$new.Range("E10").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E10").PasteSpecial(-4122) | Out-Null
$new.Range("F10").Value = "In Progress"
$style27.Copy() | Out-Null
$new.Range("F10").PasteSpecial(-4122) | Out-Null
$new.Range("G10").Value = "Still need to alter to show compared magic when looking at magic section"
$style14.Copy() | Out-Null
$new.Range("G10").PasteSpecial(-4122) | Out-Null
$new.Range("H10").Value = "Done"
$style28.Copy() | Out-Null
$new.Range("H10").PasteSpecial(-4122) | Out-Null

# Row 11
$new.Range("A11").Value = "Add iconselect boolean to prevent multiple icons being selected"
$style14.Copy() | Out-Null
$new.Range("A11").PasteSpecial(-4122) | Out-Null
$new.Range("B11").Value = 0.25
$style14.Copy() | Out-Null
$new.Range("B11").PasteSpecial(-4122) | Out-Null
$new.Range("C11").Value = 0.25
$style14.Copy() | Out-Null
$new.Range("C11").PasteSpecial(-4122) | Out-Null
$new.Range("D11").Value = 0
$style14.Copy() | Out-Null
$new.Range("D11").PasteSpecial(-4122) | Out-Null
$new.Range("E11").Value = "James"
$style14.Copy() | Out-Null
$new.Range("E11").PasteSpecial(-4122) | Out-Null
$new.Range("F11").Value = "Done"
$style28.Copy() | Out-Null
$new.Range("F11").PasteSpecial(-4122) | Out-Null

# Row 12
$style14.Copy() | Out-Null
$new.Range("A12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("B12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("C12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("E12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("F12").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("G12").PasteSpecial(-4122) | Out-Null

# Row 13
$new.Range("A13").Value = "Total Hours Assigned"
$style14.Copy() | Out-Null
$new.Range("A13").PasteSpecial(-4122) | Out-Null
$new.Range("B13").Formula = "=SUM(B3:B10)"
$style14.Copy() | Out-Null
$new.Range("B13").PasteSpecial(-4122) | Out-Null
$new.Range("C13").Formula = "=SUM(C3:C12)"
$style14.Copy() | Out-Null
$new.Range("C13").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("D13").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("E13").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("F13").PasteSpecial(-4122) | Out-Null
$style14.Copy() | Out-Null
$new.Range("G13").PasteSpecial(-4122) | Out-Null

# ---- Sheet view / selection ------------------------------------------
$new.Range("C6").Select()

# ------------------------------------------------------------------
# 2. Update "Tasks 01-28 to 02-04": add the new "Over/Under" column D
#    with computed values, and change the active selection.
# ------------------------------------------------------------------
$old.Range("D2").Value = -0.5
$old.Range("D3").Value = -0.5
$old.Range("D4").Value = 0
$old.Range("D5").Value = -1
$old.Range("D6").Value = 0
$old.Range("D7").Value = -0.5
$old.Range("D8").Value = 0
$old.Range("D9").Value = -0.75
$old.Range("D10").Value = -0.25
$old.Range("D11").Value = -0.25
$old.Range("D12").Value = -0.25
$old.Range("D13").Value = -0.75
$old.Range("D14").Value = 0

$old.Range("M9:M11").Select()
$old.Application.ActiveWindow.ScrollRow = 7
$old.Application.ActiveWindow.ScrollColumn = 2

